$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '330.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.47%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.26%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.701'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.12%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08017'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.55%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.483'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.44%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.605'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.77%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.935'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.89%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.922'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.87%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9213'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.72%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1245'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.61%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1933'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.40%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.682'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '16.92%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09244'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.79%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03566'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.14%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '9.63%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001296'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.70%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006312'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.25%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.366'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.09%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3456'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.41%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1373'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.67%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2697'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.42%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04439'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001258'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.37%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004462'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.80%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '1.31%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02554'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.51%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05523'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '5.59%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007522'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.20%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009923'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '10.90%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.00%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002113'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.34%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01168'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '21.87%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006802'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.17%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.44%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003068'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '6.96%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002285'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.64%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.44%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.44%'
